$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Paragraph 3 currently reads:
#   "Allows users to visualize and share data (sensors values, media, web links etc.)"
# Split it into three runs.
$para = $tr.Paragraphs(3)
$run1 = $para.Runs(1)
$run1.Text = "Allows users to visualize and "
$run2 = $run1.InsertAfter("anonymously share data and digital content ")
$run3 = $run2.InsertAfter("(sensors values, media, web links etc.)")
